# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the naive-forecaster QoQ error triangle on Sheet1 to reflect
# the corrected error values and the additional diagonal of newly
# available quarters (rows 24-52, columns B:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = [ordered]@{
    "K24" = -4.140068527788461;
    "J25" = -4.143502000649867;
    "K25" = 0.7890573502134768;
    "I26" = -4.148050779251975;
    "J26" = 0.7845085716113679;
    "K26" = 1.78685138562497;
    "H27" = -4.152012823709551;
    "I27" = 0.7805465271537924;
    "J27" = 1.782889341167394;
    "K27" = -1.836796412571299;
    "G28" = -4.190095585793672;
    "H28" = 0.7424637650696713;
    "I28" = 1.744806579083273;
    "J28" = -1.87487917465542;
    "K28" = -1.555883901694814;
    "F29" = -4.221995783645582;
    "G29" = 0.710563567217761;
    "H29" = 1.712906381231363;
    "I29" = -1.90677937250733;
    "J29" = -1.587784099546724;
    "K29" = 0.9930956379036073;
    "E30" = -4.304043343248933;
    "F30" = 0.6285160076144095;
    "G30" = 1.630858821628011;
    "H30" = -1.988826932110682;
    "I30" = -1.669831659150076;
    "J30" = 0.9110480783002558;
    "K30" = -1.413911231055208;
    "D31" = -4.361331021144436;
    "E31" = 0.5712283297189071;
    "F31" = 1.573571143732509;
    "G31" = -2.046114610006184;
    "H31" = -1.727119337045578;
    "I31" = 0.8537604004047534;
    "J31" = -1.47119890895071;
    "K31" = -0.4159398956771542;
    "C32" = -4.557181482025568;
    "D32" = 0.3753778688377755;
    "E32" = 1.377720682851377;
    "F32" = -2.241965070887316;
    "G32" = -1.92296979792671;
    "H32" = 0.6579099395236218;
    "I32" = -1.667049369831842;
    "J32" = -0.611790356558286;
    "K32" = -0.3935002837000762;
    "B33" = -5.278233986760878;
    "C33" = -0.3456746358975342;
    "D33" = 0.6566681781160675;
    "E33" = -2.963017575622626;
    "F33" = -2.64402230266202;
    "G33" = -0.06314256521168793;
    "H33" = -2.388101874567151;
    "I33" = -1.332842861293596;
    "J33" = -1.114552788435386;
    "K33" = -0.02381999081273611;
    "B34" = 0.1209180521064689;
    "C34" = 1.12326086612007;
    "D34" = -2.496424887618622;
    "E34" = -2.177429614658016;
    "F34" = 0.4034501227923152;
    "G34" = -1.921509186563148;
    "H34" = -0.8662501732895925;
    "I34" = -0.6479601004313827;
    "J34" = 0.442772697191267;
    "K34" = 0.8594960153664678;
    "B35" = 1.017201315193215;
    "C35" = -2.602484438545478;
    "D35" = -2.283489165584871;
    "E35" = 0.29739057186546;
    "F35" = -2.027568737490003;
    "G35" = -0.9723097242164478;
    "H35" = -0.754019651358238;
    "I35" = 0.3367131462644118;
    "J35" = 0.7534364644396125;
    "K35" = -0.8313463317717158;
    "B36" = -2.767155053718898;
    "C36" = -2.448159780758292;
    "D36" = 0.1327199566920393;
    "E36" = -2.192239352663424;
    "F36" = -1.136980339389869;
    "G36" = -0.9186902665316586;
    "H36" = 0.1720425310909911;
    "I36" = 0.5887658492661919;
    "J36" = -0.9960169469451364;
    "K36" = -0.04879898327473942;
    "B37" = -2.239654652406992;
    "C37" = 0.341225085043339;
    "D37" = -1.983734224312125;
    "E37" = -0.9284752110385688;
    "F37" = -0.710185138180359;
    "G37" = 0.3805476594422908;
    "H37" = 0.7972709776174915;
    "I37" = -0.7875118185938368;
    "J37" = 0.1597061450765603;
    "K37" = 0.08001204887913171;
    "B38" = 0.574277577822088;
    "C38" = -1.750681731533375;
    "D38" = -0.6954227182598197;
    "E38" = -0.4771326454016099;
    "F38" = 0.6136001522210398;
    "G38" = 1.030323470396241;
    "H38" = -0.5544593258150877;
    "I38" = 0.3927586378553093;
    "J38" = 0.3130645416578807;
    "K38" = 0.02451986800082353;
    "B39" = -1.870043355038259;
    "C39" = -0.8147843417647035;
    "D39" = -0.5964942689064937;
    "E39" = 0.494238528716156;
    "F39" = 0.9109618468913568;
    "G39" = -0.6738209493199715;
    "H39" = 0.2733970143504255;
    "I39" = 0.193702918152997;
    "J39" = -0.09484175550406024;
    "K39" = -0.1043830804051818;
    "B40" = -0.5707751523983828;
    "C40" = -0.3524850795401731;
    "D40" = 0.7382477180824767;
    "E40" = 1.154971036257677;
    "F40" = -0.4298117599536508;
    "G40" = 0.5174062037167462;
    "H40" = 0.4377121075193176;
    "I40" = 0.1491674338622604;
    "J40" = 0.1396261089611388;
    "K40" = -0.8797011412147235;
    "B41" = -0.3465490807841718;
    "C41" = 0.7441837168384779;
    "D41" = 1.160907035013679;
    "E41" = -0.4238757611976496;
    "F41" = 0.5233422024727474;
    "G41" = 0.4436481062753188;
    "H41" = 0.1551034326182616;
    "I41" = 0.14556210771714;
    "J41" = -0.8737651424587223;
    "K41" = -0.8519563250702333;
    "B42" = 0.8662498299576697;
    "C42" = 1.282973148132871;
    "D42" = -0.3018096480784578;
    "E42" = 0.6454083155919392;
    "F42" = 0.5657142193945106;
    "G42" = 0.2771695457374534;
    "H42" = 0.2676282208363318;
    "I42" = -0.7516990293395305;
    "J42" = -0.7298902119510415;
    "K42" = 0.8090800184609777;
    "B43" = 1.963967031984447;
    "C43" = 0.3791842357731184;
    "D43" = 1.326402199443515;
    "E43" = 1.246708103246087;
    "F43" = 0.9581634295890297;
    "G43" = 0.9486221046879081;
    "H43" = -0.07070514548795426;
    "I43" = -0.04889632809946531;
    "J43" = 1.490073902312554;
    "K43" = 0.5380735621110603;
    "B44" = -0.2239305827602892;
    "C44" = 0.7232873809101079;
    "D44" = 0.6435932847126793;
    "E44" = 0.3550486110556221;
    "F44" = 0.3455072861545005;
    "G44" = -0.6738199640213618;
    "H44" = -0.6520111466328729;
    "I44" = 0.8869590837791463;
    "J44" = -0.06504125642234732;
    "B45" = 0.7958897816072067;
    "C45" = 0.7161956854097782;
    "D45" = 0.4276510117527209;
    "E45" = 0.4181096868515993;
    "F45" = -0.601217563324263;
    "G45" = -0.5794087459357741;
    "H45" = 0.9595614844762451;
    "I45" = 0.007561144274751519;
    "B46" = 0.7786460178510652;
    "C46" = 0.490101344194008;
    "D46" = 0.4805600192928864;
    "E46" = -0.5387672308829758;
    "F46" = -0.5169584134944869;
    "G46" = 1.022011816917532;
    "H46" = 0.07001147671603863;
    "B47" = 0.2008399526861098;
    "C47" = 0.1912986277849882;
    "D47" = -0.8280286223908742;
    "E47" = -0.8062198050023852;
    "F47" = 0.732750425409634;
    "G47" = -0.2192499147918596;
    "B48" = -0.02682703805890724;
    "C48" = -1.04615428823477;
    "D48" = -1.024345470846281;
    "E48" = 0.5146247595657385;
    "F48" = -0.437375580635755;
    "B49" = -0.9928286439961705;
    "C49" = -0.9710198266076816;
    "D49" = 0.5679504038043376;
    "E49" = -0.384049936397156;
    "B50" = -0.970113206870677;
    "C50" = 0.5688570235413422;
    "D50" = -0.3831433166601514;
    "B51" = 0.6286205760598051;
    "C51" = -0.3233797641416885;
    "B52" = -0.4742879390039576
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
